$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2302218, 0, 0, 0, 0, 5, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}

$ws.Range("F2").Select()
